# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Leve profit sheets
# (scheduled data-refresh style edit across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 125021.125
$ws.Range("I11").Value = 125021.125
$ws.Range("K11").Value = 125021.125
$ws.Range("M11").Value = -124881.125
$ws.Range("H41").Value = 861.4666999999999
$ws.Range("I41").Value = 1565.5
$ws.Range("J41").Value = 509.45
$ws.Range("K41").Value = 1565.5
$ws.Range("L41").Value = 509.45
$ws.Range("M41").Value = -1125.5
$ws.Range("N41").Value = -1389.45
$ws.Range("H51").Value = 8854.888999999999
$ws.Range("I51").Value = 11935.363
$ws.Range("J51").Value = 4014.1428
$ws.Range("K51").Value = 11935.363
$ws.Range("L51").Value = 4014.1428
$ws.Range("M51").Value = -11451.363
$ws.Range("N51").Value = -4982.1428
$ws.Range("H106").Value = 3007.2856
$ws.Range("I106").Value = 4233.1665
$ws.Range("J106").Value = 2087.875
$ws.Range("K106").Value = 4233.1665
$ws.Range("L106").Value = 2087.875
$ws.Range("M106").Value = -3602.1665
$ws.Range("N106").Value = -3349.875
$ws.Range("H111").Value = 6259256
$ws.Range("I111").Value = 10647.615
$ws.Range("J111").Value = 33336560
$ws.Range("K111").Value = 31942.845
$ws.Range("L111").Value = 100009680
$ws.Range("M111").Value = -28875.845
$ws.Range("N111").Value = -100015814
$ws.Range("H129").Value = 2405.8406
$ws.Range("J129").Value = 1209.5435
$ws.Range("L129").Value = 3628.6305
$ws.Range("N129").Value = -13628.6305
$ws.Range("H132").Value = 5439761
$ws.Range("I132").Value = 6415982
$ws.Range("J132").Value = 814.7143
$ws.Range("K132").Value = 19247946
$ws.Range("L132").Value = 2444.1429
$ws.Range("M132").Value = -19245416
$ws.Range("N132").Value = -7504.1429
$ws.Range("H135").Value = 904.3158
$ws.Range("J135").Value = 2388.1
$ws.Range("L135").Value = 21492.9
$ws.Range("N135").Value = -26562.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20724.666
$ws.Range("I32").Value = 4473.9287
$ws.Range("K32").Value = 4473.9287
$ws.Range("M32").Value = -4186.9287
$ws.Range("H74").Value = 710.1111
$ws.Range("I74").Value = 580.4
$ws.Range("J74").Value = 872.25
$ws.Range("K74").Value = 580.4
$ws.Range("L74").Value = 872.25
$ws.Range("M74").Value = 293.6
$ws.Range("N74").Value = -2620.25
$ws.Range("H77").Value = 710.1111
$ws.Range("I77").Value = 580.4
$ws.Range("J77").Value = 872.25
$ws.Range("K77").Value = 2902
$ws.Range("L77").Value = 4361.25
$ws.Range("M77").Value = 1466
$ws.Range("N77").Value = -13097.25
$ws.Range("H80").Value = 27361.2
$ws.Range("J80").Value = 27361.2
$ws.Range("L80").Value = 27361.2
$ws.Range("N80").Value = -29357.2
$ws.Range("H83").Value = 27361.2
$ws.Range("J83").Value = 27361.2
$ws.Range("L83").Value = 82083.60000000001
$ws.Range("N83").Value = -92067.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 73901.78999999999
$ws.Range("I105").Value = 45891.26
$ws.Range("J105").Value = 202750.2
$ws.Range("K105").Value = 45891.26
$ws.Range("L105").Value = 202750.2
$ws.Range("M105").Value = -44144.26
$ws.Range("N105").Value = -206244.2
$ws.Range("H132").Value = 64999.855
$ws.Range("J132").Value = 64999.855
$ws.Range("L132").Value = 64999.855
$ws.Range("N132").Value = -75119.85500000001
$ws.Range("H134").Value = 3175.484
$ws.Range("I134").Value = 2808.5715
$ws.Range("J134").Value = 6600
$ws.Range("K134").Value = 8425.7145
$ws.Range("L134").Value = 19800
$ws.Range("M134").Value = -5890.7145
$ws.Range("N134").Value = -24870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1258.5306
$ws.Range("I58").Value = 1097.425
$ws.Range("J58").Value = 1974.5555
$ws.Range("K58").Value = 1097.425
$ws.Range("L58").Value = 1974.5555
$ws.Range("M58").Value = -894.425
$ws.Range("N58").Value = -2380.5555
$ws.Range("H68").Value = 17846.87
$ws.Range("J68").Value = 17846.87
$ws.Range("L68").Value = 17846.87
$ws.Range("N68").Value = -19344.87
$ws.Range("H71").Value = 17846.87
$ws.Range("J71").Value = 17846.87
$ws.Range("L71").Value = 53540.61
$ws.Range("N71").Value = -61028.61
$ws.Range("H134").Value = 1476.2632
$ws.Range("I134").Value = 1408.3334
$ws.Range("J134").Value = 1592.7142
$ws.Range("K134").Value = 4225.0002
$ws.Range("L134").Value = 4778.142599999999
$ws.Range("M134").Value = -1690.0002
$ws.Range("N134").Value = -9848.142599999999
$ws.Range("H136").Value = 1258.5306
$ws.Range("I136").Value = 1097.425
$ws.Range("J136").Value = 1974.5555
$ws.Range("K136").Value = 3292.275
$ws.Range("L136").Value = 5923.666499999999
$ws.Range("M136").Value = -742.2749999999996
$ws.Range("N136").Value = -11023.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1024.45
$ws.Range("I5").Value = 566.0952
$ws.Range("K5").Value = 1698.2856
$ws.Range("M5").Value = -1586.2856
$ws.Range("H18").Value = 461.5
$ws.Range("I18").Value = 461.6154
$ws.Range("J18").Value = 461
$ws.Range("K18").Value = 1384.8462
$ws.Range("L18").Value = 1383
$ws.Range("M18").Value = -1215.8462
$ws.Range("N18").Value = -1721
$ws.Range("H74").Value = 4100
$ws.Range("H77").Value = 4100
$ws.Range("H131").Value = 7278.1846
$ws.Range("I131").Value = 1858
$ws.Range("J131").Value = 7460.8877
$ws.Range("K131").Value = 5574
$ws.Range("L131").Value = 22382.6631
$ws.Range("M131").Value = -534
$ws.Range("N131").Value = -32462.6631
$ws.Range("H135").Value = 1024.45
$ws.Range("I135").Value = 566.0952
$ws.Range("K135").Value = 5094.8568
$ws.Range("M135").Value = -2559.8568
$ws.Range("H140").Value = 4804.613
$ws.Range("I140").Value = 6228.2
$ws.Range("J140").Value = 2216.2727
$ws.Range("K140").Value = 18684.6
$ws.Range("L140").Value = 6648.8181
$ws.Range("M140").Value = -13504.6
$ws.Range("N140").Value = -17008.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 47669870
$ws.Range("I80").Value = 91004150
$ws.Range("J80").Value = 2162.7
$ws.Range("K80").Value = 91004150
$ws.Range("L80").Value = 2162.7
$ws.Range("M80").Value = -91003152
$ws.Range("N80").Value = -4158.7
$ws.Range("H83").Value = 47669870
$ws.Range("I83").Value = 91004150
$ws.Range("J83").Value = 2162.7
$ws.Range("K83").Value = 455020750
$ws.Range("L83").Value = 10813.5
$ws.Range("M83").Value = -455015758
$ws.Range("N83").Value = -20797.5
$ws.Range("H132").Value = 2435.8027
$ws.Range("I132").Value = 2246.5283
$ws.Range("J132").Value = 2993.111
$ws.Range("K132").Value = 6739.5849
$ws.Range("L132").Value = 8979.332999999999
$ws.Range("M132").Value = -4209.5849
$ws.Range("N132").Value = -14039.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3916.2693
$ws.Range("I132").Value = 4046.25
$ws.Range("K132").Value = 12138.75
$ws.Range("M132").Value = -9608.75
$ws.Range("H136").Value = 1056.1569
$ws.Range("I136").Value = 854.57446
$ws.Range("J136").Value = 3424.75
$ws.Range("K136").Value = 2563.72338
$ws.Range("L136").Value = 10274.25
$ws.Range("M136").Value = -13.72338000000036
$ws.Range("N136").Value = -15374.25
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1962.7142
$ws.Range("I122").Value = 1126.3334
$ws.Range("J122").Value = 2590
$ws.Range("K122").Value = 3379.0002
$ws.Range("L122").Value = 7770
$ws.Range("M122").Value = -929.0001999999999
$ws.Range("N122").Value = -12670
$ws.Range("H132").Value = 2428
$ws.Range("I132").Value = 2304.1538
$ws.Range("J132").Value = 3233
$ws.Range("K132").Value = 6912.4614
$ws.Range("L132").Value = 9699
$ws.Range("M132").Value = -4382.4614
$ws.Range("N132").Value = -14759
